# Fix spatial relationship bug: add a lookup/verification helper table
# (COCO class id -> raw label -> trimmed label) in columns M:O of the
# "Sheet2" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# --- Rows 5-14: id / raw label (with TRIM formula in column O) ---
$ws.Range("M5").Value = 13
$ws.Range("N5").Value = " bench"
$ws.Range("O5").Formula = "=TRIM(N5)"

$ws.Range("M6").Value = 56
$ws.Range("N6").Value = " chair"

$ws.Range("M7").Value = 57
$ws.Range("N7").Value = " couch"

$ws.Range("M8").Value = 58
$ws.Range("N8").Value = " potted plant"

$ws.Range("M9").Value = 59
$ws.Range("N9").Value = " bed"

$ws.Range("M10").Value = 60
$ws.Range("N10").Value = " dining table"

$ws.Range("M11").Value = 62
$ws.Range("N11").Value = " tv"

$ws.Range("M12").Value = 68
$ws.Range("N12").Value = " microwave"

$ws.Range("M13").Value = 69
$ws.Range("N13").Value = " oven"

$ws.Range("M14").Value = 72
$ws.Range("N14").Value = " refrigerator"

# Shared TRIM formula covering O6:O14
$ws.Range("O6:O14").Formula = "=TRIM(N6)"

# --- Rows 21-46: id / raw label only (no formula column) ---
$ws.Range("M21").Value = 15
$ws.Range("N21").Value = " cat"

$ws.Range("M22").Value = 16
$ws.Range("N22").Value = " dog"

$ws.Range("M23").Value = 24
$ws.Range("N23").Value = " backpack"

$ws.Range("M24").Value = 25
$ws.Range("N24").Value = " umbrella"

$ws.Range("M25").Value = 26
$ws.Range("N25").Value = " handbag"

$ws.Range("M26").Value = 27
$ws.Range("N26").Value = " tie"

$ws.Range("M27").Value = 28
$ws.Range("N27").Value = " suitcase"

$ws.Range("M28").Value = 32
$ws.Range("N28").Value = " sports ball"

$ws.Range("M29").Value = 39
$ws.Range("N29").Value = " bottle"

$ws.Range("M30").Value = 41
$ws.Range("N30").Value = " cup"

$ws.Range("M31").Value = 42
$ws.Range("N31").Value = " fork"

$ws.Range("M32").Value = 43
$ws.Range("N32").Value = " knife"

$ws.Range("M33").Value = 44
$ws.Range("N33").Value = " spoon"

$ws.Range("M34").Value = 45
$ws.Range("N34").Value = " bowl"

$ws.Range("M35").Value = 63
$ws.Range("N35").Value = " laptop"

$ws.Range("M36").Value = 64
$ws.Range("N36").Value = " mouse"

$ws.Range("M37").Value = 65
$ws.Range("N37").Value = " remote"

$ws.Range("M38").Value = 66
$ws.Range("N38").Value = " keyboard"

$ws.Range("M39").Value = 67
$ws.Range("N39").Value = " cell phone"

$ws.Range("M40").Value = 73
$ws.Range("N40").Value = " book"

$ws.Range("M41").Value = 74
$ws.Range("N41").Value = " clock"

$ws.Range("M42").Value = 75
$ws.Range("N42").Value = " vase"

$ws.Range("M43").Value = 76
$ws.Range("N43").Value = " scissors"

$ws.Range("M44").Value = 77
$ws.Range("N44").Value = " teddy bear"

$ws.Range("M45").Value = 78
$ws.Range("N45").Value = " hair drier"

$ws.Range("M46").Value = 79
$ws.Range("N46").Value = " toothbrush"

# --- View state: zoom out and move selection to the new helper table ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 85
$ws.Range("M41:N46").Select()
